# Group 3 plot commit
# Divide every value in C2:C33 by 5 (cfu_count_undiluted column) and
# move the active selection from B1 to F2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    1100000,
    1440000,
    1060000,
    700000,
    7600000,
    7800000,
    7000000,
    8600000,
    2600000,
    3400000,
    2400000,
    3200000,
    138000,
    154000,
    142000,
    162000,
    94000,
    104000,
    112000,
    84000,
    82000,
    102000,
    86000,
    82000,
    920000,
    560000,
    700000,
    800000,
    740000,
    1060000,
    1040000,
    800000
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $newValues[$i]
}

$ws.Range("F2").Select()
